$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Planning")
$cells = @("F11","F15","F19","F24","F25","F28","F32","F35","F36","F37","F38","F39","F40","F41","F42","F43","F44","F45")
foreach ($addr in $cells) {
    $ws2.Range($addr).Value = "FG/HG"
}
